$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# New text used by this revision
# ------------------------------------------------------------------
$childCasePath = "src/test/resources/runCase/child-case.xlsx"
$newJson = "{`n    `"target`":`"src/test/resources/runCase/child-case.xlsx`",`n    `"cases`":`"1`",`n    `"type`":`"xlsx`"`n}"

# ------------------------------------------------------------------
# Column D: give it the same width as column C, and make sure its
# border/shading formatting (thin box border) matches column C before
# we touch alignment, by lifting the formats from the matching C cell
# in each row.
# ------------------------------------------------------------------
$ws.Range("D1").EntireColumn.ColumnWidth = 59.75

$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("C2").Copy()
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("C2").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("D3").WrapText = $true

# ------------------------------------------------------------------
# Cell values
# ------------------------------------------------------------------
# Row 1 (header row) - new D1 repeats the "runCase" header
$ws.Range("D1").Value = $ws.Range("C1").Value()

# Row 2 (setProperty / json row) - C2 now points at the child workbook
$ws.Range("C2").Value = $childCasePath
$ws.Range("D2").Value = ""

# Row 3 (TestCase row) - C3 becomes the case count, D3 the run config
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = $newJson

# ------------------------------------------------------------------
# Alignment: every bordered cell (header + data cells) is now left
# aligned horizontally, instead of the old mix of centered / default.
# (This runtime only honours the first area of a comma-separated
# Range, so each contiguous block is aligned with its own call.)
# ------------------------------------------------------------------
$ws.Range("A1:D1").HorizontalAlignment = -4131
$ws.Range("A2:B2").HorizontalAlignment = -4131
$ws.Range("A3:C3").HorizontalAlignment = -4131
$ws.Range("C2:D2").HorizontalAlignment = -4131
$ws.Range("D3:D3").HorizontalAlignment = -4131

# ------------------------------------------------------------------
# Row heights: the "tall" wrap-text row moves from row 2 to row 3.
# ------------------------------------------------------------------
$ws.Rows.Item(2).AutoFit()
$ws.Range("A3").EntireRow.RowHeight = 93.75

# ------------------------------------------------------------------
# Selection moves to the newly added D3 cell
# ------------------------------------------------------------------
$ws.Range("D3").Select()
